$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing accuracy value for Lasso Regression+normalization+ lag1+PCA(2)
$ws.Range("C7").Value = 83.519033481304604

# Add new data preprocessing row
$ws.Range("B17").Value = '`'

# Update selected cell
[void]$ws.Range("E8").Select()
